$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.476685643196106
$ws.Range("B1").Value = 1.683289051055908
$ws.Range("C1").Value = 2.094408750534058
$ws.Range("D1").Value = 2.155177116394043
$ws.Range("E1").Value = 1.439373731613159
